$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Add a date value in C1 (2021-04-21), formatted as a date.
$c1 = $ws.Range("C1")
$c1.NumberFormat = "mm-dd-yy"
$c1.Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
